# Atualizado por script em 05-11-2023 14:45
# Appends the new match row (row 97) to Sheet1, matching the formatting
# of the preceding data row (row 96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (bold/bordered index style) and column E (date/time number format)
# carry explicit cell styles in this sheet; copy just those formats down from
# the last existing data row before writing the new values.
$ws.Range("A96").Copy()
$ws.Range("A97").PasteSpecial(-4122)
$ws.Range("E96").Copy()
$ws.Range("E97").PasteSpecial(-4122)

$ws.Range("A97").Value = 96
$ws.Range("B97").Value = "serbia"
$ws.Range("C97").Value = "super-liga"
$ws.Range("D97").Value = "2023-2024"
$ws.Range("E97").Value = 45235.625
$ws.Range("F97").Value = "Vojvodina"
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = "Zeleznicar Pancevo"
$ws.Range("I97").Value = 2
$ws.Range("J97").Value = 1.42
$ws.Range("K97").Value = "02/11/2023 08:13"
$ws.Range("L97").Value = 1.37
$ws.Range("M97").Value = "05/11/2023 14:52"
$ws.Range("N97").Value = 4.26
$ws.Range("O97").Value = "02/11/2023 08:13"
$ws.Range("P97").Value = 4.39
$ws.Range("Q97").Value = "05/11/2023 14:52"
$ws.Range("R97").Value = 5.66
$ws.Range("S97").Value = "02/11/2023 08:13"
$ws.Range("T97").Value = 8.79
$ws.Range("U97").Value = "05/11/2023 14:52"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/serbia/super-liga/vojvodina-zeleznicar-pancevo/rahtDksK/"
